# Auto-generated market-price refresh for Leve crafting profit sheets
# (mirrors the scheduled runner that repopulates currentAveragePrice* / Leve* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 40016
$ws.Range("I11").Value = 40016
$ws.Range("K11").Value = 40016
$ws.Range("M11").Value = -39876

$ws.Range("H41").Value = 888.7059
$ws.Range("I41").Value = 226.25
$ws.Range("K41").Value = 226.25
$ws.Range("M41").Value = 213.75

$ws.Range("H43").Value = 33344328
$ws.Range("J43").Value = 17583
$ws.Range("L43").Value = 17583
$ws.Range("N43").Value = -17721

$ws.Range("H51").Value = 62629564
$ws.Range("I51").Value = 204999.4
$ws.Range("J51").Value = 166670500
$ws.Range("K51").Value = 204999.4
$ws.Range("L51").Value = 166670500
$ws.Range("M51").Value = -204515.4
$ws.Range("N51").Value = -166671468

$ws.Range("H86").Value = 266669570
$ws.Range("I86").Value = 333335650
$ws.Range("K86").Value = 333335650
$ws.Range("M86").Value = -333334527

$ws.Range("H89").Value = 266669570
$ws.Range("I89").Value = 333335650
$ws.Range("K89").Value = 1666678250
$ws.Range("M89").Value = -1666672634

$ws.Range("H113").Value = 2683.1667
$ws.Range("I113").Value = 1619.8
$ws.Range("K113").Value = 1619.8
$ws.Range("M113").Value = 1634.2

$ws.Range("H131").Value = 7164.857
$ws.Range("I131").Value = 7164.857
$ws.Range("K131").Value = 21494.571
$ws.Range("M131").Value = -16454.571

$ws.Range("H132").Value = 1804.1428
$ws.Range("I132").Value = 1587.7391
$ws.Range("K132").Value = 4763.2173
$ws.Range("M132").Value = -2233.2173

$ws.Range("H135").Value = 1794.5555
$ws.Range("I135").Value = 1608.3636
$ws.Range("K135").Value = 14475.2724
$ws.Range("M135").Value = -11940.2724

$ws.Range("H137").Value = 1375968.9
$ws.Range("I137").Value = 4398.8706
$ws.Range("J137").Value = 5274115
$ws.Range("K137").Value = 13196.6118
$ws.Range("L137").Value = 15822345
$ws.Range("M137").Value = -10646.6118
$ws.Range("N137").Value = -15827445

$ws.Range("H138").Value = 6437.2285
$ws.Range("J138").Value = 4014.7144
$ws.Range("L138").Value = 12044.1432
$ws.Range("N138").Value = -22324.1432

$ws.Range("H141").Value = 3727.2
$ws.Range("I141").Value = 3703.3333
$ws.Range("J141").Value = 3822.6667
$ws.Range("K141").Value = 11109.9999
$ws.Range("L141").Value = 11468.0001
$ws.Range("M141").Value = -5929.999899999999
$ws.Range("N141").Value = -21828.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 773.3333
$ws.Range("I2").Value = 748
$ws.Range("J2").Value = 900
$ws.Range("K2").Value = 748
$ws.Range("L2").Value = 900
$ws.Range("M2").Value = -635
$ws.Range("N2").Value = -1126

$ws.Range("H32").Value = 13143.333
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 13143.333
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 13143.333
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -13717.333

$ws.Range("H45").Value = 64301.75
$ws.Range("I45").Value = 78617.53999999999
$ws.Range("J45").Value = 2266.6667
$ws.Range("K45").Value = 78617.53999999999
$ws.Range("L45").Value = 2266.6667
$ws.Range("M45").Value = -78240.53999999999
$ws.Range("N45").Value = -3020.6667

$ws.Range("I61").Value = 24666.275
$ws.Range("J61").Value = 3994834.5
$ws.Range("K61").Value = 24666.275
$ws.Range("L61").Value = 3994834.5
$ws.Range("M61").Value = -24454.275
$ws.Range("N61").Value = -3995258.5

$ws.Range("H109").Value = 49950
$ws.Range("J109").Value = 49950
$ws.Range("L109").Value = 49950
$ws.Range("N109").Value = -52724

$ws.Range("H116").Value = 773.3333
$ws.Range("I116").Value = 748
$ws.Range("J116").Value = 900
$ws.Range("K116").Value = 748
$ws.Range("L116").Value = 900
$ws.Range("M116").Value = 1546
$ws.Range("N116").Value = -5488

$ws.Range("H122").Value = 2694.125
$ws.Range("I122").Value = 2694.125
$ws.Range("K122").Value = 8082.375
$ws.Range("M122").Value = -5632.375

$ws.Range("I136").Value = 24666.275
$ws.Range("J136").Value = 3994834.5
$ws.Range("K136").Value = 73998.82500000001
$ws.Range("L136").Value = 11984503.5
$ws.Range("M136").Value = -71448.82500000001
$ws.Range("N136").Value = -11989603.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 773.3333
$ws.Range("I3").Value = 748
$ws.Range("J3").Value = 900
$ws.Range("K3").Value = 748
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = -634
$ws.Range("N3").Value = -1128

$ws.Range("J80").Value = 3454.25
$ws.Range("L80").Value = 3454.25
$ws.Range("N80").Value = -5450.25

$ws.Range("J83").Value = 3454.25
$ws.Range("L83").Value = 17271.25
$ws.Range("N83").Value = -27255.25

$ws.Range("H86").Value = 7237.25
$ws.Range("I86").Value = 2350
$ws.Range("J86").Value = 12124.5
$ws.Range("K86").Value = 2350
$ws.Range("L86").Value = 12124.5
$ws.Range("M86").Value = -1227
$ws.Range("N86").Value = -14370.5

$ws.Range("H89").Value = 7237.25
$ws.Range("I89").Value = 2350
$ws.Range("J89").Value = 12124.5
$ws.Range("K89").Value = 11750
$ws.Range("L89").Value = 60622.5
$ws.Range("M89").Value = -6134
$ws.Range("N89").Value = -71854.5

$ws.Range("H105").Value = 9654.200000000001
$ws.Range("I105").Value = 6377.6787
$ws.Range("K105").Value = 6377.6787
$ws.Range("M105").Value = -4630.6787

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5604.8335
$ws.Range("J122").Value = 2600
$ws.Range("L122").Value = 7800
$ws.Range("N122").Value = -12700

$ws.Range("H132").Value = 2983.2222
$ws.Range("I132").Value = 2985.7144
$ws.Range("K132").Value = 8957.143199999999
$ws.Range("M132").Value = -6427.143199999999

$ws.Range("H134").Value = 2717.087
$ws.Range("I134").Value = 2564.5
$ws.Range("K134").Value = 7693.5
$ws.Range("M134").Value = -5158.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 280.53845
$ws.Range("I2").Value = 209.83333
$ws.Range("J2").Value = 341.14285
$ws.Range("K2").Value = 1258.99998
$ws.Range("L2").Value = 2046.8571
$ws.Range("M2").Value = -1145.99998
$ws.Range("N2").Value = -2272.8571

$ws.Range("H7").Value = 278
$ws.Range("I7").Value = 278
$ws.Range("K7").Value = 834
$ws.Range("M7").Value = -722

$ws.Range("H80").Value = 2500
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 2500
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H97").Value = 1384.8
$ws.Range("I97").Value = 360.5
$ws.Range("J97").Value = 2067.6667
$ws.Range("K97").Value = 1081.5
$ws.Range("L97").Value = 6203.000100000001
$ws.Range("M97").Value = -585.5
$ws.Range("N97").Value = -7195.000100000001

$ws.Range("H108").Value = 200003360
$ws.Range("J108").Value = 5000
$ws.Range("L108").Value = 15000
$ws.Range("N108").Value = -20760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 38463816
$ws.Range("I102").Value = 41668800
$ws.Range("K102").Value = 41668800
$ws.Range("M102").Value = -41667178

$ws.Range("H106").Value = 34593
$ws.Range("J106").Value = 34593
$ws.Range("L106").Value = 34593
$ws.Range("N106").Value = -37117

$ws.Range("H126").Value = 2502.8333
$ws.Range("I126").Value = 2470.4443
$ws.Range("K126").Value = 7411.3329
$ws.Range("M126").Value = -4941.3329

$ws.Range("H132").Value = 755810.0600000001
$ws.Range("I132").Value = 1177
$ws.Range("K132").Value = 3531
$ws.Range("M132").Value = -1001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 22857.123
$ws.Range("I136").Value = 36312.414
$ws.Range("J136").Value = 3346.95
$ws.Range("K136").Value = 108937.242
$ws.Range("L136").Value = 10040.85
$ws.Range("M136").Value = -106387.242
$ws.Range("N136").Value = -15140.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 31467.8
$ws.Range("J104").Value = 31467.8
$ws.Range("L104").Value = 31467.8
$ws.Range("N104").Value = -38455.8

$ws.Range("H107").Value = 2858281
$ws.Range("I107").Value = 985.375
$ws.Range("K107").Value = 2956.125
$ws.Range("M107").Value = -1036.125

$ws.Range("H113").Value = 1033.5454
$ws.Range("I113").Value = 920
$ws.Range("K113").Value = 2760
$ws.Range("M113").Value = -590

$ws.Range("H132").Value = 19609492
$ws.Range("I132").Value = 29412696
$ws.Range("J132").Value = 3087.7058
$ws.Range("K132").Value = 88238088
$ws.Range("L132").Value = 9263.117400000001
$ws.Range("M132").Value = -88235558
$ws.Range("N132").Value = -14323.1174

$ws.Range("H135").Value = 90363.39999999999
$ws.Range("J135").Value = 90363.39999999999
$ws.Range("L135").Value = 90363.39999999999
$ws.Range("N135").Value = -100503.4

$ws.Range("H136").Value = 704.75
$ws.Range("I136").Value = 559.74286
$ws.Range("J136").Value = 1719.8
$ws.Range("K136").Value = 1679.22858
$ws.Range("L136").Value = 5159.4
$ws.Range("M136").Value = 870.77142
$ws.Range("N136").Value = -10259.4
